$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell values scraped from the updated cryptos.xlsx sheet (columns B-E, rows 2-51).
# Using an explicit text NumberFormat while assigning the value (and clearing the
# formatting afterwards) keeps numeric-looking strings such as "40.91" or "0.0691"
# stored as text, matching the original inline-string cell type instead of letting
# Excel auto-convert them to numbers.
$updates = [ordered]@{
    'D2' = '34.967.35'
    'E2' = '  +0.28%  '
    'D3' = '1.845.91'
    'E3' = '  +2.09%  '
    'E4' = '  +0.17%  '
    'D5' = '232.67'
    'E5' = '  +0.54%  '
    'E6' = '  +2.29%  '
    'E7' = '  +0.22%  '
    'D8' = '40.91'
    'E8' = '  +4.77%  '
    'D9' = '0.330'
    'E9' = '  +3.93%  '
    'D10' = '0.0691'
    'E10' = '  +1.97%  '
    'D11' = '0.0986'
    'E11' = '  -0.81%  '
    'D12' = '2.115.07'
    'E12' = '  +2.22%  '
    'E13' = '  +4.91%  '
    'D14' = '1.853.44'
    'E14' = '  +2.43%  '
    'D15' = '0.673'
    'E15' = '  +2.16%  '
    'D16' = '4.66'
    'E16' = '  +2.32%  '
    'D17' = '34.997.37'
    'E17' = '  +0.49%  '
    'D18' = '70.05'
    'E18' = '  +1.60%  '
    'E19' = '  +1.62%  '
    'D20' = '240.48'
    'D21' = '12.22'
    'E21' = '  +3.95%  '
    'E22' = '  +2.72%  '
    'E23' = '  +0.23%  '
    'E24' = '  +1.47%  '
    'D25' = '172.84'
    'E25' = '  +0.53%  '
    'D26' = '7.85'
    'E26' = '  +1.36%  '
    'D27' = '17.49'
    'E27' = '  +1.79%  '
    'E28' = '  +4.02%  '
    'D29' = '1.62'
    'E29' = '  +4.74%  '
    'E30' = '  +0.16%  '
    'E31' = '  +1.35%  '
    'E32' = '  +0.01%  '
    'E33' = '  +1.32%  '
    'D34' = '1.61'
    'E34' = '  +23.74%  '
    'E35' = '  +11.54%  '
    'D36' = '0.748'
    'E36' = '  +10.12%  '
    'E37' = '  +6.77%  '
    'D38' = '1.07'
    'E38' = '  +12.09%  '
    'D39' = '90.07'
    'E39' = '  -0.87%  '
    'D40' = '1.349.66'
    'E40' = '  +3.38%  '
    'D41' = '0.0196'
    'E41' = '  +3.02%  '
    'D42' = '14.78'
    'E42' = '  +2.66%  '
    'E43' = '  +3.54%  '
    'E44' = '  -1.58%  '
    'E45' = '  +2.58%  '
    'E46' = '  +4.04%  '
    'D47' = '6.32'
    'E47' = '  +3.14%  '
    'D48' = '2.032.60'
    'E48' = '  +2.15%  '
    'B49' = 'THORChain'
    'C49' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D49' = '3.41'
    'E49' = '  +18.75%  '
    'B50' = 'PaxDollar'
    'C50' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D50' = '1.01'
    'E50' = '  +0.26%  '
    'D51' = '0.0671'
    'E51' = '  -0.39%  '

}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
